$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Status moved from "In Translation" to "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Latest HO Xliff Generate Date refreshed
$wsOverview.Range("G2").Value = "2016-08-18 11:01:48"
# Widen the Status columns (E & F) to fit "Ready for handoff"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-18 11:01:42"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-18 11:01:48"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
